$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Text
$text = $text.Replace("✅ 1000 Bs = 9.01 = 37297.63 pesos", "✅ 1000 Bs = 8.98 = 37204.23 pesos")
$text = $text.Replace("✅ 37297.63 pesos = 8.96 = 954.62 Bs", "✅ 37204.23 pesos = 9.0 = 960.8 Bs")
$cell.Value = $text

# --- Sheet "tasas": update the rate table values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 111.318
$ws2.Range("O10").Value = 4141.5
$ws2.Range("N12").Value = 4135
$ws2.Range("O12").Value = 106.786
